$wb = $excel.ActiveWorkbook

$ttd = $wb.Worksheets.Item("TTD")
$yahoo = $wb.Worksheets.Item("Yahoo")

# Remember the currency number format currently used on G3 (row that will become row 4)
$currencyFormat = $ttd.Cells.Item(3, 7).NumberFormat

# --- Update TTD (sheet2) rows 3 and 4 ---
# New row 4 becomes what used to be row 3, with the Buyer Member ID changed to a plain number
$ttd.Cells.Item(4, 1).Value = 31804
$ttd.Cells.Item(4, 2).Value = 31804
$ttd.Cells.Item(4, 3).Value = "Premium Range Purchasers"
$ttd.Cells.Item(4, 4).Value = "Users who prefer premium branded ranges over supermarket own brand"
$ttd.Cells.Item(4, 5).Value = "Buyable"
$ttd.Cells.Item(4, 6).Value = 3
$ttd.Cells.Item(4, 7).NumberFormat = $currencyFormat
$ttd.Cells.Item(4, 7).Value = 1.5
$ttd.Cells.Item(4, 8).Value = "UK Kantar Media TGI > Grocery Shopping > Premium Range Purchasers"

# New row 3 becomes what used to be row 4, with the Parent Segment ID changed to "taxoapitest"
$ttd.Cells.Item(3, 1).Value = 32048
$ttd.Cells.Item(3, 2).Value = "taxoapitest"
$ttd.Cells.Item(3, 3).Value = "TV Channels Watched Live (Last 30 Days)"
$ttd.Cells.Item(3, 4).ClearContents()
$ttd.Cells.Item(3, 5).Value = "Not Buyable"
$ttd.Cells.Item(3, 6).Value = 3
$ttd.Cells.Item(3, 7).ClearFormats()
$ttd.Cells.Item(3, 7).Value = 0
$ttd.Cells.Item(3, 8).Value = "Media > TV And Film > TV Channels Watched Live (Last 30 Days)"

# --- Add a new "Segment Description" column to Yahoo (sheet3) ---
$yahoo.Cells.Item(1, 3).Value = "Segment Description"
$yahoo.Cells.Item(2, 3).Value = "Required"
$yahoo.Cells.Item(3, 3).Value = "Just random segment 4444"
$yahoo.Cells.Item(4, 3).Value = "Just random segment 1111"
$yahoo.Cells.Item(5, 3).Value = "Just random segment 2222"
$yahoo.Cells.Item(6, 3).Value = "Just random segment 3333"

$yahoo.Range("C1").Copy()
$yahoo.Cells.Item(1, 3).PasteSpecial(-4122) | Out-Null
$yahoo.Cells.Item(1, 3).Style = $yahoo.Cells.Item(1, 2).Style
$yahoo.Cells.Item(2, 3).Style = $yahoo.Cells.Item(2, 2).Style

# --- Update sheet selections / active tab ---
# Yahoo is no longer the selected tab, but keep its own remembered selection at C7
$yahoo.Activate()
$yahoo.Range("C7").Select()

# TTD becomes the active / selected tab with H9 selected
$ttd.Activate()
$ttd.Range("H9").Select()
